$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Alpna sheet (sheet2): fill in row 120 (Alpna / MT2012013 / Nothing)
# ---------------------------------------------------------------------
$alpna = $wb.Worksheets.Item("Alpna")
$alpna.Range("B120").Value = "Alpna"
$alpna.Range("C120").Value = "MT2012013"
$alpna.Range("E120").Value = "Nothing"

# ---------------------------------------------------------------------
# 2) Sravani sheet (sheet7): populate the timesheet for team member 7
# ---------------------------------------------------------------------
$sravani = $wb.Worksheets.Item("Sravani")

$sravani.Range("A1").Value = "Team member 7: "
$sravani.Range("A2").Value = "Name"
$sravani.Range("B2").Value = "RollNo"
$sravani.Range("C2").Value = "Date"
$sravani.Range("D2").Value = "Phase"
$sravani.Range("E2").Value = "Activity"
$sravani.Range("F2").Value = "Time Spent(Hrs)"

$name = "Sai Naga Sravani Peraka"
$roll = "MT2012122"

$data = @(
  @(40912, "Requirements for Ebay", "Going through Ebay site and looked at some of its features by creating an user`naccount", 1),
  @(40913, "Requirements for Ebay", "Studied some literature about Ebay", 2),
  @(40914, "Nothing", "Nothing", "-"),
  @(40915, "Requirements for Ebay", "Studied some literature about Ebay", 1),
  @(40916, "Requirements for Ebay", 'Preparing the document "Know your friend" along with Team mates', 2),
  @(40917, "Nothing", "Nothing", "-"),
  @(40918, "Nothing", "Nothing", "-"),
  @(40919, "Nothing", "Nothing", "-"),
  @(40920, "Requirements for Ebay", "Understanding Requirements", 2),
  @(40921, "Nothing", "Nothing", "-"),
  @(40922, "Nothing", "Nothing", "-"),
  @(40923, "Nothing", "Nothing", "-"),
  @(40924, "Nothing", "Nothing", "-"),
  @(40925, "Requirements for Ebay", "Understanding Requirements and thought of some feasible usecases`n by discussing with team members", 3),
  @(40926, "Nothing", "Nothing", "-"),
  @(40927, "Nothing", "Nothing", "-"),
  @(40928, "Nothing", "Nothing", "-"),
  @(40929, "Requirements for Ebay", "Preparing Usecase diagrams", 3),
  @(40930, "Nothing", "Nothing", "-"),
  @(40931, "Nothing", "Nothing", "-"),
  @(40932, "Requirements for Ebay", "Preparing Interaction Stories", 4),
  @(40933, "Nothing", "Nothing", "-"),
  @(40934, "Nothing", "Nothing", "-"),
  @(40935, "Requirements for Ebay", "Preparing Flow of events", 3),
  @(40936, "Requirements for Ebay", "Preparing Test Transactions", 3),
  @(40937, "Requirements for MiniProject", "Understanding Requirements for Student Profile Project", 2),
  @(40938, "Database Design for Miniproject", "Brainstorming for Database tables according to requirements", 1),
  @(41305, "Formal Meeting(with Sir)", 1.5, 3),
  @(40940, "Database Design for Miniproject", "Finalizing the Database schema by discussing with Teammembers and mentor", 2),
  @(40941, "SVN Repository Setup", "Created the SVN repositories for both Ebay and StudentProfile and done the`ninitial SVN commit by creating the dynamic web projects and basic folder structures", 3)
)

$row = 3
foreach ($entry in $data) {
    $sravani.Range("A$row").Value = $name
    $sravani.Range("B$row").Value = $roll
    $sravani.Range("C$row").Value = $entry[0]
    $sravani.Range("D$row").Value = $entry[1]
    $sravani.Range("E$row").Value = $entry[2]
    $sravani.Range("F$row").Value = $entry[3]
    $row = $row + 1
}

$sravani.Range("E3").Select()

# ---------------------------------------------------------------------
# 3) Switch the active tab from Deepthi to Alpna, select F120 there
# ---------------------------------------------------------------------
$alpna.Activate()
$alpna.Range("F120").Select()
